# Mise à jour de l'application
# Adds a new attendance date column (BB) for 2025-09-26 right after the
# existing last date column (BA, 2025-09-24), fills in the attendance
# status per player for that date, and updates the active selection.
#
# NOTE: values are written first (so the COUNTA/COUNTIF summary formulas
# in columns B-J recalculate against the new data), and formatting is
# copied from column BA afterwards so the new cells end up with the same
# style as the rest of the date columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) New date header for the added column (2025-09-26).
$ws.Range("BB1").Value = 45926

# 2) Attendance values for the new date, one per player row.
#    (Row 12's player has no record for this date, same as column BA.)
$ws.Range("BB2").Value = "P"
$ws.Range("BB3").Value = "P"
$ws.Range("BB4").Value = "P"
$ws.Range("BB5").Value = "RH"
$ws.Range("BB6").Value = "M"
$ws.Range("BB7").Value = "P"
$ws.Range("BB8").Value = "P"
$ws.Range("BB9").Value = "P"
$ws.Range("BB10").Value = "RH"
$ws.Range("BB11").Value = "B"
$ws.Range("BB13").Value = "B"
$ws.Range("BB14").Value = "P"
$ws.Range("BB15").Value = "P"
$ws.Range("BB16").Value = "P"
$ws.Range("BB17").Value = "P"
$ws.Range("BB18").Value = "P"
$ws.Range("BB19").Value = "P"
$ws.Range("BB20").Value = "P"
$ws.Range("BB21").Value = "B"
$ws.Range("BB22").Value = "P"
$ws.Range("BB23").Value = "B"
$ws.Range("BB24").Value = "P"
$ws.Range("BB25").Value = "B"
$ws.Range("BB26").Value = "P"
$ws.Range("BB27").Value = "REP"
$ws.Range("BB28").Value = "P"
$ws.Range("BB29").Value = "M"

# 3) Copy the formatting of column BA onto the new column BB.
$ws.Range("BA1:BA29").Copy()
$ws.Range("BB1:BB29").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 12 has no entry in column BA (that player's data stops at column AX),
# so after the style copy there must be no BB12 cell either.
$ws.Range("BB12").Clear()

# 4) Update the active cell selection as left by the author after editing.
$ws.Range("BD24").Select()
